$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits: sample/borehole IDs renumbered -------------------------
# ST-06 / ST-07 (test-boring style IDs) renamed to SS-06 / SS-07
$ws.Range("A9").Value = "SS-06"
$ws.Range("A10").Value = "SS-07"

# Second "SS-12" entry (row 16) becomes "SS-13"
$ws.Range("A16").Value = "SS-13"

# DB-19..DB-22 (bedrock samples) renamed to SS-19..SS-22
$ws.Range("A22").Value = "SS-19"
$ws.Range("A23").Value = "SS-20"
$ws.Range("A24").Value = "SS-21"
$ws.Range("A25").Value = "SS-22"

# --- View changes ------------------------------------------------------
# Reset scroll position back to the top-left corner (removes topLeftCell="A3")
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
# Zoom 259% -> 150%
$excel.ActiveWindow.Zoom = 150
# Selection moves from A15 to A6
$ws.Range("A6").Select()
